$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 100048.836
$ws.Range("I40").Value = 501150
$ws.Range("J40").Value = 3784.56
$ws.Range("K40").Value = 501150
$ws.Range("L40").Value = 3784.56
$ws.Range("M40").Value = -500975
$ws.Range("N40").Value = -4134.559999999999
$ws.Range("H92").Value = 41667664
$ws.Range("I92").Value = 50001180
$ws.Range("J92").Value = 98.5
$ws.Range("K92").Value = 50001180
$ws.Range("L92").Value = 98.5
$ws.Range("M92").Value = -49999932
$ws.Range("N92").Value = -2594.5
$ws.Range("H116").Value = 55579440
$ws.Range("I116").Value = 71457224
$ws.Range("K116").Value = 71457224
$ws.Range("M116").Value = -71453782
$ws.Range("H125").Value = 2373.1428
$ws.Range("I125").Value = 1381.8572
$ws.Range("J125").Value = 2703.5715
$ws.Range("K125").Value = 12436.7148
$ws.Range("L125").Value = 24332.1435
$ws.Range("M125").Value = -9976.7148
$ws.Range("N125").Value = -29252.1435
$ws.Range("H132").Value = 3588.1636
$ws.Range("I132").Value = 3612.6416
$ws.Range("J132").Value = 2939.5
$ws.Range("K132").Value = 10837.9248
$ws.Range("L132").Value = 8818.5
$ws.Range("M132").Value = -8307.924800000001
$ws.Range("N132").Value = -13878.5
$ws.Range("H135").Value = 1858.4584
$ws.Range("I135").Value = 1658.6666
$ws.Range("J135").Value = 3257
$ws.Range("K135").Value = 14927.9994
$ws.Range("L135").Value = 29313
$ws.Range("M135").Value = -12392.9994
$ws.Range("N135").Value = -34383
$ws.Range("H137").Value = 17527.02
$ws.Range("I137").Value = 19513.457
$ws.Range("J137").Value = 2297.6667
$ws.Range("K137").Value = 58540.371
$ws.Range("L137").Value = 6893.000100000001
$ws.Range("M137").Value = -55990.371
$ws.Range("N137").Value = -11993.0001
$ws.Range("H138").Value = 1832.44
$ws.Range("I138").Value = 1041.9546
$ws.Range("J138").Value = 2453.5356
$ws.Range("K138").Value = 3125.8638
$ws.Range("L138").Value = 7360.6068
$ws.Range("M138").Value = 2014.1362
$ws.Range("N138").Value = -17640.6068
$ws.Range("H141").Value = 1181.1
$ws.Range("I141").Value = 1130.7059
$ws.Range("J141").Value = 1466.6666
$ws.Range("K141").Value = 3392.1177
$ws.Range("L141").Value = 4399.9998
$ws.Range("M141").Value = 1787.8823
$ws.Range("N141").Value = -14759.9998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18349314
$ws.Range("J32").Value = 17862516
$ws.Range("L32").Value = 17862516
$ws.Range("N32").Value = -17863090
$ws.Range("H41").Value = 6396
$ws.Range("I41").Value = 6396
$ws.Range("K41").Value = 6396
$ws.Range("M41").Value = -5982
$ws.Range("H63").Value = 3563.4546
$ws.Range("I63").Value = 2485.4285
$ws.Range("K63").Value = 2485.4285
$ws.Range("M63").Value = -1799.4285
$ws.Range("H66").Value = 3563.4546
$ws.Range("I66").Value = 2485.4285
$ws.Range("K66").Value = 12427.1425
$ws.Range("M66").Value = -8995.1425
$ws.Range("H74").Value = 2257.8333
$ws.Range("I74").Value = 2361.4119
$ws.Range("J74").Value = 1817.625
$ws.Range("K74").Value = 2361.4119
$ws.Range("L74").Value = 1817.625
$ws.Range("M74").Value = -1487.4119
$ws.Range("N74").Value = -3565.625
$ws.Range("H77").Value = 2257.8333
$ws.Range("I77").Value = 2361.4119
$ws.Range("J77").Value = 1817.625
$ws.Range("K77").Value = 11807.0595
$ws.Range("L77").Value = 9088.125
$ws.Range("M77").Value = -7439.059499999999
$ws.Range("N77").Value = -17824.125
$ws.Range("H122").Value = 3790.842
$ws.Range("J122").Value = 5779.2354
$ws.Range("L122").Value = 17337.7062
$ws.Range("N122").Value = -22237.7062
$ws.Range("H132").Value = 2322.0217
$ws.Range("I132").Value = 1602.1765
$ws.Range("K132").Value = 4806.529500000001
$ws.Range("M132").Value = -2276.529500000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 29939.5
$ws.Range("J109").Value = 29939.5
$ws.Range("L109").Value = 29939.5
$ws.Range("N109").Value = -32713.5
$ws.Range("H134").Value = 1135849.8
$ws.Range("I134").Value = 1300406.6
$ws.Range("K134").Value = 3901219.8
$ws.Range("M134").Value = -3898684.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 4039.3333
$ws.Range("I5").Value = 120
$ws.Range("K5").Value = 120
$ws.Range("M5").Value = -8
$ws.Range("H15").Value = 878
$ws.Range("J15").Value = 750
$ws.Range("L15").Value = 750
$ws.Range("N15").Value = -1090
$ws.Range("H31").Value = 3810.8728
$ws.Range("I31").Value = 1662.081
$ws.Range("J31").Value = 8227.833000000001
$ws.Range("K31").Value = 1662.081
$ws.Range("L31").Value = 8227.833000000001
$ws.Range("M31").Value = -1367.081
$ws.Range("N31").Value = -8817.833000000001
$ws.Range("H34").Value = 3810.8728
$ws.Range("I34").Value = 1662.081
$ws.Range("J34").Value = 8227.833000000001
$ws.Range("K34").Value = 1662.081
$ws.Range("L34").Value = 8227.833000000001
$ws.Range("M34").Value = -1460.081
$ws.Range("N34").Value = -8631.833000000001
$ws.Range("H99").Value = 999
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 999
$ws.Range("M99").Value = 499
$ws.Range("H122").Value = 4003075.5
$ws.Range("I122").Value = 5003247
$ws.Range("K122").Value = 15009741
$ws.Range("M122").Value = -15007291
$ws.Range("H126").Value = 999
$ws.Range("I126").Value = 999
$ws.Range("K126").Value = 2997
$ws.Range("M126").Value = -527
$ws.Range("H132").Value = 2648.7932
$ws.Range("I132").Value = 2597.5588
$ws.Range("J132").Value = 2832.158
$ws.Range("K132").Value = 7792.676399999999
$ws.Range("L132").Value = 8496.474
$ws.Range("M132").Value = -5262.676399999999
$ws.Range("N132").Value = -13556.474

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3149.1667
$ws.Range("I58").Value = 1779
$ws.Range("K58").Value = 5337
$ws.Range("M58").Value = -5209
$ws.Range("H64").Value = 2120.25
$ws.Range("I64").Value = 2120.25
$ws.Range("K64").Value = 6360.75
$ws.Range("M64").Value = -6090.75
$ws.Range("H67").Value = 2120.25
$ws.Range("I67").Value = 2120.25
$ws.Range("K67").Value = 6360.75
$ws.Range("M67").Value = -5424.75
$ws.Range("H140").Value = 2632.4644
$ws.Range("I140").Value = 2077.3333
$ws.Range("J140").Value = 3273
$ws.Range("K140").Value = 6231.999899999999
$ws.Range("L140").Value = 9819
$ws.Range("M140").Value = -1051.999899999999
$ws.Range("N140").Value = -20179
$ws.Range("H141").Value = 3738.2222
$ws.Range("I141").Value = 2731.3333
$ws.Range("K141").Value = 8193.999899999999
$ws.Range("M141").Value = -3013.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1939.25
$ws.Range("I122").Value = 1102.5
$ws.Range("J122").Value = 2441.3
$ws.Range("K122").Value = 3307.5
$ws.Range("L122").Value = 7323.900000000001
$ws.Range("M122").Value = -857.5
$ws.Range("N122").Value = -12223.9
$ws.Range("H132").Value = 4961.1626
$ws.Range("I132").Value = 4737
$ws.Range("K132").Value = 14211
$ws.Range("M132").Value = -11681

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13895957
$ws.Range("I40").Value = 16672949
$ws.Range("K40").Value = 16672949
$ws.Range("M40").Value = -16672813
$ws.Range("H45").Value = 30041
$ws.Range("I45").Value = 30041
$ws.Range("K45").Value = 30041
$ws.Range("M45").Value = -29634
$ws.Range("H122").Value = 18989.53
$ws.Range("I122").Value = 23425.77
$ws.Range("K122").Value = 70277.31
$ws.Range("M122").Value = -67827.31
$ws.Range("H132").Value = 91450.52
$ws.Range("I132").Value = 109247.21
$ws.Range("K132").Value = 327741.63
$ws.Range("M132").Value = -325211.63
$ws.Range("H136").Value = 3281.1667
$ws.Range("I136").Value = 2671.875
$ws.Range("J136").Value = 4499.75
$ws.Range("K136").Value = 8015.625
$ws.Range("L136").Value = 13499.25
$ws.Range("M136").Value = -5465.625
$ws.Range("N136").Value = -18599.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 90028
$ws.Range("I49").Value = 30056
$ws.Range("J49").Value = 150000
$ws.Range("K49").Value = 30056
$ws.Range("L49").Value = 150000
$ws.Range("M49").Value = -29826
$ws.Range("N49").Value = -150460
$ws.Range("H122").Value = 52635732
$ws.Range("I122").Value = 66669300
$ws.Range("J122").Value = 9843.5
$ws.Range("K122").Value = 200007900
$ws.Range("L122").Value = 29530.5
$ws.Range("M122").Value = -200005450
$ws.Range("N122").Value = -34430.5
$ws.Range("H132").Value = 2904
$ws.Range("I132").Value = 2624.5667
$ws.Range("K132").Value = 7873.7001
$ws.Range("M132").Value = -5343.7001
$ws.Range("H136").Value = 65884.31
$ws.Range("I136").Value = 3044.3333
$ws.Range("K136").Value = 9132.999899999999
$ws.Range("M136").Value = -6582.999899999999
